$wb = $excel.ActiveWorkbook

# --- Section_A (sheet1) ---
$wsA = $wb.Worksheets.Item("Section_A")

$wsA.Range("B2").Value = "Free"
$wsA.Range("C2").Value = "CS307"
$wsA.Range("D2").Value = "Free"
$wsA.Range("E2").Value = "Free"
$wsA.Range("F2").Value = "DS303"

$wsA.Range("B3").Value = "CS307"
$wsA.Range("C3").Value = "Free"
$wsA.Range("D3").Value = "DS303"
$wsA.Range("E3").Value = "DS303"
$wsA.Range("F3").Value = "Free"

$wsA.Range("B5").Value = "Free"
$wsA.Range("C5").Value = "Free"
$wsA.Range("D5").Value = "DS302"
$wsA.Range("E5").Value = "DS302"
$wsA.Range("F5").Value = "DS302"

$wsA.Range("B6").Value = "CS307 (Tutorial)"
$wsA.Range("C6").Value = "DS302 (Tutorial)"
$wsA.Range("D6").Value = "Free"
$wsA.Range("E6").Value = "Free"
$wsA.Range("F6").Value = "Free"

$wsA.Range("B7").Value = "Free"
$wsA.Range("C7").Value = "Free"
$wsA.Range("D7").Value = "Free"
$wsA.Range("E7").Value = "Free"
$wsA.Range("F7").Value = "CS307"

$wsA.Range("B8").Value = "DS303 (Tutorial)"
$wsA.Range("C8").Value = "Free"
$wsA.Range("D8").Value = "Free"
$wsA.Range("E8").Value = "Free"
$wsA.Range("F8").Value = "Free"

# --- Section_B (sheet2) ---
$wsB = $wb.Worksheets.Item("Section_B")

$wsB.Range("B2").Value = "DS303"
$wsB.Range("C2").Value = "DS303"
$wsB.Range("D2").Value = "CS307"
$wsB.Range("E2").Value = "DS303"
$wsB.Range("F2").Value = "Free"

$wsB.Range("B3").Value = "Free"
$wsB.Range("C3").Value = "CS307"
$wsB.Range("D3").Value = "Free"
$wsB.Range("E3").Value = "Free"
$wsB.Range("F3").Value = "DS302"

$wsB.Range("B5").Value = "Free"
$wsB.Range("C5").Value = "Free"
$wsB.Range("D5").Value = "DS302"
$wsB.Range("E5").Value = "Free"
$wsB.Range("F5").Value = "Free"

$wsB.Range("B6").Value = "CS307 (Tutorial)"
$wsB.Range("C6").Value = "Free"
$wsB.Range("D6").Value = "DS303 (Tutorial)"
$wsB.Range("E6").Value = "Free"
$wsB.Range("F6").Value = "Free"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "DS302"
$wsB.Range("D7").Value = "Free"
$wsB.Range("E7").Value = "CS307"
$wsB.Range("F7").Value = "Free"

$wsB.Range("B8").Value = "Free"
$wsB.Range("C8").Value = "Free"
$wsB.Range("D8").Value = "Free"
$wsB.Range("E8").Value = "DS302 (Tutorial)"
$wsB.Range("F8").Value = "Free"
